# v2.2 update according to new comments
#
# 1. LH_TC_USERHOME sheet: the shared "Test Data" text (email + password) used
#    across the test-case rows gets its password corrected to include the
#    "@" symbol.
# 2. VESRION HISTORY sheet: the v2.1 row's "Updated section" note is corrected,
#    and a new v2.2 row is appended documenting this very update.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: LH_TC_USERHOME ---------------------------------------------
$ws1 = $wb.Worksheets.Item("LH_TC_USERHOME")

$oldTestData = "E-mail: user1@example.com`nPassword: CorrectPassword123"
$newTestData = "E-mail: user1@example.com`nPassword: CorrectPassword@123"

for ($r = 9; $r -le 17; $r++) {
    $cell = $ws1.Cells.Item($r, 6)
    if ($cell.Text -eq $oldTestData) {
        $cell.Value = $newTestData
    }
}

# --- Sheet 2: VESRION HISTORY ----------------------------------------------
$ws2 = $wb.Worksheets.Item("VESRION HISTORY")

# Append the new v2.2 row, matching row 3's formatting.
$ws2.Rows.Item(3).Copy()
$ws2.Rows.Item(4).PasteSpecial(-4104)

$ws2.Cells.Item(4, 1).Value = "v2.2"
$ws2.Cells.Item(4, 2).Value = "Hala Eldaly"

$updatedNote = "update according to new comments"

# Row 3 (the v2.1 entry) gets its "Updated section" note corrected.
$ws2.Cells.Item(3, 3).Value = $updatedNote
$ws2.Cells.Item(4, 3).Value = $updatedNote
$ws2.Cells.Item(4, 4).Value = (Get-Date -Year 2025 -Month 5 -Day 14 -Hour 0 -Minute 0 -Second 0)

$ws2.Range("D5").Select()
